# Append a new scrape run (2025-11-10 01:23:14 JST) to the "ランサーズ" sheet.
# Two fresh listings are inserted at the top of the data (rows 10-11); the
# stalest listing that fell out of the top-N ("小売店向けシステム性能試験")
# is dropped, and the remaining previously-tracked rows shift down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$timestamp = "2025-11-10 01:23:14"

# Widen column B by one character (51 -> 52), matching the COM
# ColumnWidth<->stored-width offset measured on this sheet (~0.83).
$colB = $ws.Columns.Item(2)
$colB.ColumnWidth = $colB.ColumnWidth + 1

# Insert a single blank row at row 10: this pushes the old row 10
# ("小売店向けシステム性能試験", which is being replaced/dropped) down to
# row 11, and everything below it down by one as well -- exactly matching
# the target layout where old rows 11-15 become new rows 12-16.
$ws.Rows.Item(10).Insert()

# --- New row 10: 完全在宅GASエンジニア募集... -----------------------------
$ws.Cells.Item(10, 1).Value2 = $timestamp
$ws.Cells.Item(10, 2).Value2 = "完全在宅GASエンジニア募集/課題テストからご依頼/時給1,163円~業務フロー効率化をお任せします"
$ws.Cells.Item(10, 3).Value2 = "システム開発"
$ws.Cells.Item(10, 4).Value2 = "~ 5,000 円 / 固定"
$ws.Cells.Item(10, 5).Value2 = "期限情報なし"
$ws.Cells.Item(10, 6).Value2 = "https://www.lancers.jp/work/detail/5416665"
$ws.Cells.Item(10, 7).Value2 = 70
$ws.Cells.Item(10, 8).Value2 = "◆効率化"

# --- New row 11: 【業務効率化】SlackとHubSpotの活用支援... (overwrites the
#     shifted-down old row 10 data, which is being replaced entirely) -----
$ws.Cells.Item(11, 1).Value2 = $timestamp
$ws.Cells.Item(11, 2).Value2 = "【業務効率化】SlackとHubSpotの活用支援をお願いします"
$ws.Cells.Item(11, 3).Value2 = "システム開発"
$ws.Cells.Item(11, 4).Value2 = "5,000 円 ~ 10,000 円 / 固定"
$ws.Cells.Item(11, 5).Value2 = "期限情報なし"
$ws.Cells.Item(11, 6).Value2 = "https://www.lancers.jp/work/detail/5430436"
$ws.Cells.Item(11, 7).Value2 = 70
$ws.Cells.Item(11, 8).Value2 = "◆効率化"

# --- Refresh the A column timestamp on every already-tracked row ---------
# (rows 2-9 are untouched pre-existing listings; rows 12-16 after the shift
# are the old rows 11-15). Every row in this scrape run shares one stamp.
foreach ($r in 2..9 + 12..16) {
    $ws.Cells.Item($r, 1).Value2 = $timestamp
}

# --- Rebuild hyperlinks for column F (rows 2-16). The row Insert() above
#     shifts cell contents correctly but leaves the saved <hyperlinks> table
#     pointing at the pre-insert row positions, so every entry needs to be
#     re-created against the new layout. -----------------------------------
$ws.Hyperlinks.Delete()
for ($r = 2; $r -le 16; $r++) {
    $target = $ws.Cells.Item($r, 6).Value2
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $target)
    $ws.Cells.Item($r, 6).Value2 = $target
    # Hyperlinks.Add() stamps a freshly-minted "Hyperlink" style; re-apply
    # the named style so column F keeps using the workbook's existing one
    # instead of accumulating a duplicate.
    $ws.Cells.Item($r, 6).Style = "Hyperlink"
}

Write-Output "done"
